$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 2570
$ws.Cells.Item(28, 10).Value = 5999.5
$ws.Cells.Item(28, 12).Value = 5999.5
$ws.Cells.Item(28, 14).Value = -6969.5
$ws.Cells.Item(43, 8).Value = 5066.3335
$ws.Cells.Item(43, 10).Value = 5066.3335
$ws.Cells.Item(43, 12).Value = 5066.3335
$ws.Cells.Item(43, 14).Value = -5204.3335
$ws.Cells.Item(103, 8).Value = 748.4
$ws.Cells.Item(103, 9).Value = 716.6667
$ws.Cells.Item(103, 10).Value = 796
$ws.Cells.Item(103, 11).Value = 2150.0001
$ws.Cells.Item(103, 12).Value = 2388
$ws.Cells.Item(103, 13).Value = -1564.0001
$ws.Cells.Item(103, 14).Value = -3560
$ws.Cells.Item(116, 8).Value = 3449.5
$ws.Cells.Item(116, 9).Value = 3150
$ws.Cells.Item(116, 10).Value = 3749
$ws.Cells.Item(116, 11).Value = 3150
$ws.Cells.Item(116, 12).Value = 3749
$ws.Cells.Item(116, 13).Value = 292
$ws.Cells.Item(116, 14).Value = -10633
$ws.Cells.Item(141, 8).Value = 8317.6
$ws.Cells.Item(141, 9).Value = 8317.6
$ws.Cells.Item(141, 11).Value = 24952.8
$ws.Cells.Item(141, 13).Value = -19772.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13458.454
$ws.Cells.Item(32, 9).Value = 13304.75
$ws.Cells.Item(32, 10).Value = 14995.5
$ws.Cells.Item(32, 11).Value = 13304.75
$ws.Cells.Item(32, 12).Value = 14995.5
$ws.Cells.Item(32, 13).Value = -13017.75
$ws.Cells.Item(32, 14).Value = -15569.5
$ws.Cells.Item(45, 8).Value = 5946.625
$ws.Cells.Item(45, 9).Value = 5946.625
$ws.Cells.Item(45, 11).Value = 5946.625
$ws.Cells.Item(45, 13).Value = -5569.625
$ws.Cells.Item(98, 8).Value = 20521
$ws.Cells.Item(98, 10).Value = 20521
$ws.Cells.Item(98, 12).Value = 20521
$ws.Cells.Item(98, 14).Value = -26511
$ws.Cells.Item(139, 8).Value = 0
$ws.Cells.Item(139, 10).Value = 0
$ws.Cells.Item(139, 12).Value = 0
$ws.Cells.Item(139, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2336.5
$ws.Cells.Item(20, 9).Value = 1898.8572
$ws.Cells.Item(20, 11).Value = 1898.8572
$ws.Cells.Item(20, 13).Value = -1651.8572
$ws.Cells.Item(86, 8).Value = 27299.6
$ws.Cells.Item(86, 9).Value = 38999.332
$ws.Cells.Item(86, 10).Value = 9750
$ws.Cells.Item(86, 11).Value = 38999.332
$ws.Cells.Item(86, 12).Value = 9750
$ws.Cells.Item(86, 13).Value = -37876.332
$ws.Cells.Item(86, 14).Value = -11996
$ws.Cells.Item(89, 8).Value = 27299.6
$ws.Cells.Item(89, 9).Value = 38999.332
$ws.Cells.Item(89, 10).Value = 9750
$ws.Cells.Item(89, 11).Value = 194996.66
$ws.Cells.Item(89, 12).Value = 48750
$ws.Cells.Item(89, 13).Value = -189380.66
$ws.Cells.Item(89, 14).Value = -59982
$ws.Cells.Item(105, 8).Value = 3163.5715
$ws.Cells.Item(105, 9).Value = 3124.1667
$ws.Cells.Item(105, 11).Value = 3124.1667
$ws.Cells.Item(105, 13).Value = -1377.1667
$ws.Cells.Item(122, 8).Value = 60000
$ws.Cells.Item(122, 10).Value = 60000
$ws.Cells.Item(122, 12).Value = 60000
$ws.Cells.Item(122, 14).Value = -69800
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5324.154
$ws.Cells.Item(31, 9).Value = 7741.909
$ws.Cells.Item(31, 11).Value = 7741.909
$ws.Cells.Item(31, 13).Value = -7446.909
$ws.Cells.Item(34, 8).Value = 5324.154
$ws.Cells.Item(34, 9).Value = 7741.909
$ws.Cells.Item(34, 11).Value = 7741.909
$ws.Cells.Item(34, 13).Value = -7539.909
$ws.Cells.Item(62, 8).Value = 3466.6667
$ws.Cells.Item(62, 9).Value = 3466.6667
$ws.Cells.Item(62, 11).Value = 3466.6667
$ws.Cells.Item(62, 13).Value = -2842.6667
$ws.Cells.Item(65, 8).Value = 3466.6667
$ws.Cells.Item(65, 9).Value = 3466.6667
$ws.Cells.Item(65, 11).Value = 17333.3335
$ws.Cells.Item(65, 13).Value = -14213.3335
$ws.Cells.Item(107, 8).Value = 1082.5454
$ws.Cells.Item(107, 9).Value = 1008.06665
$ws.Cells.Item(107, 11).Value = 1008.06665
$ws.Cells.Item(107, 13).Value = 911.93335
$ws.Cells.Item(132, 8).Value = 5828
$ws.Cells.Item(132, 9).Value = 5714.75
$ws.Cells.Item(132, 11).Value = 17144.25
$ws.Cells.Item(132, 13).Value = -14614.25
$ws.Cells.Item(141, 8).Value = 362981.44
$ws.Cells.Item(141, 10).Value = 362981.44
$ws.Cells.Item(141, 12).Value = 362981.44
$ws.Cells.Item(141, 14).Value = -373341.44
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 1022.2
$ws.Cells.Item(132, 9).Value = 1027.75
$ws.Cells.Item(132, 11).Value = 9249.75
$ws.Cells.Item(132, 13).Value = -6719.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(52, 8).Value = 20000
$ws.Cells.Item(52, 9).Value = 20000
$ws.Cells.Item(52, 11).Value = 20000
$ws.Cells.Item(52, 13).Value = -19741
$ws.Cells.Item(70, 8).Value = 6474
$ws.Cells.Item(70, 9).Value = 6474
$ws.Cells.Item(70, 10).Value = 0
$ws.Cells.Item(70, 11).Value = 6474
$ws.Cells.Item(70, 12).Value = 0
$ws.Cells.Item(70, 13).Value = -6204
$ws.Cells.Item(70, 14).ClearContents()
$ws.Cells.Item(73, 8).Value = 6474
$ws.Cells.Item(73, 9).Value = 6474
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 6474
$ws.Cells.Item(73, 12).Value = 0
$ws.Cells.Item(73, 13).Value = -5538
$ws.Cells.Item(73, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 2396.1
$ws.Cells.Item(102, 9).Value = 1829
$ws.Cells.Item(102, 11).Value = 1829
$ws.Cells.Item(102, 13).Value = -207
$ws.Cells.Item(105, 8).Value = 21282.143
$ws.Cells.Item(105, 10).Value = 21282.143
$ws.Cells.Item(105, 12).Value = 21282.143
$ws.Cells.Item(105, 14).Value = -28270.143
$ws.Cells.Item(141, 8).Value = 118057.4
$ws.Cells.Item(141, 10).Value = 118057.4
$ws.Cells.Item(141, 12).Value = 118057.4
$ws.Cells.Item(141, 14).Value = -128417.4
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(14, 8).Value = 8996.25
$ws.Cells.Item(14, 9).Value = 500
$ws.Cells.Item(14, 10).Value = 11828.333
$ws.Cells.Item(14, 11).Value = 500
$ws.Cells.Item(14, 12).Value = 11828.333
$ws.Cells.Item(14, 13).Value = -328
$ws.Cells.Item(14, 14).Value = -12172.333
$ws.Cells.Item(76, 8).Value = 16600
$ws.Cells.Item(76, 10).Value = 16600
$ws.Cells.Item(76, 12).Value = 16600
$ws.Cells.Item(76, 14).Value = -17276
$ws.Cells.Item(79, 8).Value = 16600
$ws.Cells.Item(79, 10).Value = 16600
$ws.Cells.Item(79, 12).Value = 16600
$ws.Cells.Item(79, 14).Value = -18940
$ws.Cells.Item(136, 8).Value = 2746
$ws.Cells.Item(136, 9).Value = 1996
$ws.Cells.Item(136, 10).Value = 4996
$ws.Cells.Item(136, 11).Value = 5988
$ws.Cells.Item(136, 12).Value = 14988
$ws.Cells.Item(136, 13).Value = -3438
$ws.Cells.Item(136, 14).Value = -20088
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 852.5357
$ws.Cells.Item(107, 10).Value = 969.1667
$ws.Cells.Item(107, 12).Value = 2907.5001
$ws.Cells.Item(107, 14).Value = -6747.5001
$ws.Cells.Item(123, 8).Value = 59992.25
$ws.Cells.Item(123, 10).Value = 59992.25
$ws.Cells.Item(123, 12).Value = 59992.25
$ws.Cells.Item(123, 14).Value = -69792.25
$ws.Cells.Item(135, 8).Value = 37858
$ws.Cells.Item(135, 10).Value = 37858
$ws.Cells.Item(135, 12).Value = 37858
$ws.Cells.Item(135, 14).Value = -47998
$ws.Cells.Item(136, 8).Value = 52134.57
$ws.Cells.Item(136, 9).Value = 51635
$ws.Cells.Item(136, 10).Value = 53966.332
$ws.Cells.Item(136, 11).Value = 154905
$ws.Cells.Item(136, 12).Value = 161898.996
$ws.Cells.Item(136, 13).Value = -152355
$ws.Cells.Item(136, 14).Value = -166998.996
